# This script updates the prediction results table on the active sheet:
# each trade's accuracy (B), start/end day range (C/D), perc_trade (I),
# n_BO (J) and n_hidd_layer (K) are refreshed with the latest run's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    # Force the cell to stay a text value (avoids Excel auto-converting
    # strings that look like numbers/percentages), then restore the
    # cell style back to Normal so no extra formatting is introduced.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2
$ws.Range("B2").Value = 0.5537604456824513
Set-TextCell "C2" "20190101"
Set-TextCell "D2" "20210130"
Set-TextCell "I2" "100%"
$ws.Range("J2").Value = 25
$ws.Range("K2").Value = 7

# Row 3
$ws.Range("B3").Value = 0.5606617647058824
Set-TextCell "C3" "20190101"
Set-TextCell "D3" "20210130"
Set-TextCell "I3" "61%"
$ws.Range("J3").Value = 25
$ws.Range("K3").Value = 3

# Row 4
$ws.Range("B4").Value = 0.5778025655326269
Set-TextCell "C4" "20190101"
Set-TextCell "D4" "20210130"
Set-TextCell "I4" "100%"
$ws.Range("J4").Value = 25
$ws.Range("K4").Value = 3

# Row 5
$ws.Range("B5").Value = 0.5887276785714286
Set-TextCell "C5" "20190101"
Set-TextCell "D5" "20210130"
Set-TextCell "I5" "100%"
$ws.Range("J5").Value = 25
$ws.Range("K5").Value = 3

# Row 6
$ws.Range("B6").Value = 0.5991066443327749
Set-TextCell "C6" "20190101"
Set-TextCell "D6" "20210130"
Set-TextCell "I6" "100%"
$ws.Range("J6").Value = 25
$ws.Range("K6").Value = 3

# Row 7
$ws.Range("B7").Value = 0.605586592178771
Set-TextCell "C7" "20190101"
Set-TextCell "D7" "20210130"
Set-TextCell "I7" "100%"
$ws.Range("J7").Value = 25
$ws.Range("K7").Value = 6

# Row 8
$ws.Range("B8").Value = 0.5975164353542732
Set-TextCell "C8" "20190101"
Set-TextCell "D8" "20210130"
Set-TextCell "I8" "77%"
$ws.Range("J8").Value = 25
$ws.Range("K8").Value = 2

# Row 9
$ws.Range("B9").Value = 0.6211832061068703
Set-TextCell "C9" "20190101"
Set-TextCell "D9" "20210130"
Set-TextCell "I9" "59%"
$ws.Range("J9").Value = 25
$ws.Range("K9").Value = 6

# Row 10
$ws.Range("B10").Value = 0.6435185185185185
Set-TextCell "C10" "20190101"
Set-TextCell "D10" "20210130"
Set-TextCell "I10" "12%"
$ws.Range("J10").Value = 25
$ws.Range("K10").Value = 2

# Row 11
$ws.Range("B11").Value = 0.6159014557670772
Set-TextCell "C11" "20190101"
Set-TextCell "D11" "20210130"
Set-TextCell "I11" "100%"
$ws.Range("J11").Value = 25
$ws.Range("K11").Value = 1

# Row 12
$ws.Range("B12").Value = 0.6179271708683474
Set-TextCell "C12" "20190101"
Set-TextCell "D12" "20210130"
Set-TextCell "I12" "100%"
$ws.Range("J12").Value = 25
$ws.Range("K12").Value = 6

# Row 13
$ws.Range("B13").Value = 0.6293706293706294
Set-TextCell "C13" "20190101"
Set-TextCell "D13" "20210130"
Set-TextCell "I13" "8%"
$ws.Range("J13").Value = 25
$ws.Range("K13").Value = 3

# Row 14
$ws.Range("B14").Value = 0.6141334828939988
Set-TextCell "C14" "20190101"
Set-TextCell "D14" "20210130"
Set-TextCell "I14" "100%"
$ws.Range("J14").Value = 25
$ws.Range("K14").Value = 4

# Row 15
$ws.Range("B15").Value = 0.5618153364632238
Set-TextCell "C15" "20190101"
Set-TextCell "D15" "20210130"
Set-TextCell "I15" "36%"
$ws.Range("J15").Value = 25
$ws.Range("K15").Value = 5

# Row 16 (keeps the older date range, n_BO stays 12)
Set-TextCell "C16" "20160101"
Set-TextCell "D16" "20191231"
Set-TextCell "I16" "5%"

# Row 17 (keeps the older date range, n_BO stays 12)
Set-TextCell "C17" "20160101"
Set-TextCell "D17" "20191231"
Set-TextCell "I17" "11%"
